$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Presentation-level slide guides (p:extLst / p15:sldGuideLst)
#    Two new guides: a horizontal one at 2160 and a vertical one at 2880.
#    (ppHorizontalGuide = 1, ppVerticalGuide = 2)
# ---------------------------------------------------------------------------
try {
    $guides = $p.Guides
    $hGuide = $guides.Add(1, 2160)
    if ($hGuide -ne $null) {
        $hGuide.Orientation = 1
        $hGuide.Position = 2160
    }
    $vGuide = $guides.Add(2, 2880)
    if ($vGuide -ne $null) {
        $vGuide.Orientation = 2
        $vGuide.Position = 2880
    }
} catch {
}

# ---------------------------------------------------------------------------
# 2. Refresh the cached "update automatically" date placeholder text
#    (the a:fld type="datetimeFigureOut" run) on the slide master and on
#    every slide layout, from 6/5/17 -> 2/17/2018.
# ---------------------------------------------------------------------------
$newDate = "2/17/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shape.HasTextFrame) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
